# Backlog sheet: a "Gemensam meny för hemsidan" task (previously row 18) is
# moved up to row 13 (marked done / "Bra"), and the tasks that used to sit in
# rows 13-17 each shift down one row to make room (rows 14-18).
#
# Column A (the running index numbers) and everything from row 19 downward
# are untouched - only B:G for rows 13-18 change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current (pre-edit) values of B:G for rows 13-18 so the shift
# below reads from a stable source instead of cells we've already written.
$rows = @(13, 14, 15, 16, 17, 18)
$cols = @("B", "C", "D", "F", "G")

$snapshot = @{}
foreach ($r in $rows) {
    foreach ($c in $cols) {
        $addr = "$c$r"
        $snapshot[$addr] = $ws.Range($addr).Value2
    }
}

# New row 13 = old row 18's content (the navbar/menu task), marked as done.
$ws.Range("B13").Value2 = $snapshot["B18"]
$ws.Range("C13").Value2 = $snapshot["C18"]
$ws.Range("D13").Value2 = $snapshot["D18"]
$ws.Range("E13").Style = "Bra"
$ws.Range("F13").Value2 = $snapshot["F13"]
$ws.Range("G13").ClearContents()

# Rows 14-18 = old rows 13-17, shifted down by one.
for ($i = 0; $i -lt 5; $i++) {
    $oldRow = 13 + $i
    $newRow = 14 + $i

    $ws.Range("B$newRow").Value2 = $snapshot["B$oldRow"]
    $ws.Range("C$newRow").Value2 = $snapshot["C$oldRow"]
    $ws.Range("D$newRow").Value2 = $snapshot["D$oldRow"]
    $ws.Range("F$newRow").Value2 = $snapshot["F$oldRow"]

    if ($snapshot["G$oldRow"]) {
        $ws.Range("G$newRow").Value2 = $snapshot["G$oldRow"]
    } else {
        $ws.Range("G$newRow").ClearContents()
    }
}

# The column holding the comments (G) was widened a bit.
$ws.Columns.Item(7).ColumnWidth = 83.5

# Last thing the user touched was cell G13.
$ws.Range("G13").Select() | Out-Null
